# Update the "想去人数" (F column) figures on the sheets that list event
# details: "展览" and "全部类型" (sheet1 and sheet4 contain identical data).

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F (rows are 1-indexed, row 1 is header)
$updates = @{
    2  = 32
    3  = 6262
    4  = 175
    6  = 38
    7  = 1886
    8  = 1417
    9  = 294
    10 = 955
    11 = 240
    12 = 5579
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
